# The edit swaps the internal "name" label recorded on the two logo
# pictures that live in the document's headers/footers:
#   - Pearson logo (footers):  name="image1.png" -> name="image2.png"
#   - BTEC logo   (headers):  name="image2.jpg" -> name="image1.jpg"
# (this only changes the wp:docPr / pic:cNvPr "name" attribute text -
#  the actual embedded media relationships are untouched)

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Update-LogoName {
    param(
        $range,
        [string]$oldName,
        [string]$newName
    )

    $xml = $range.WordOpenXML
    $target = 'name="' + $oldName + '"'
    $replacement = 'name="' + $newName + '"'
    if ($xml.IndexOf($target) -ge 0) {
        $xml = $xml.Replace($target, $replacement)
        $range.WordOpenXML = $xml
    }
}

# Footers hold the Pearson logo -> rename image1.png to image2.png
for ($i = 1; $i -le 2; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists) {
        Update-LogoName -range $ftr.Range -oldName "image1.png" -newName "image2.png"
    }
}

# Headers hold the BTEC logo -> rename image2.jpg to image1.jpg
for ($i = 1; $i -le 2; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists) {
        Update-LogoName -range $hdr.Range -oldName "image2.jpg" -newName "image1.jpg"
    }
}

Write-Host "Logo names updated"
